$wb = $excel.ActiveWorkbook

# --- unitdict sheet: add a new "DR cutoff tier 3" unit, right after
#     "DR cutoff tier 2" (row 72), by inserting a fresh row 73 and
#     shifting everything below it down by one row. ---
$ws = $wb.Worksheets.Item("unitdict")
$ws.Rows.Item(73).Insert() | Out-Null

$ws.Range("A73").Value = "DR cutoff tier 3"
$ws.Range("B73").Value = "DRcutofftier3"

# Preserve the pre-existing row-height quirk (row with ht=14.25 in a
# sea of ht=13.8 rows) which rides along with the shift: what used to
# be row 75 (14.25) is now row 76.
$ws.Rows.Item(73).RowHeight = 13.8
$ws.Rows.Item(76).RowHeight = 14.25

# Touch row 93 so the trailing blank row exists past the new last
# data row (91).
$ws.Rows.Item(93).RowHeight = 13.8

# Match the saved selection/view state on the active sheet.
$ws.Activate() | Out-Null
$ws.Range("C73").Select() | Out-Null
